$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 30
$ws.Range("K30").Value = 16
$ws.Range("L30").Value = 10.74

# Row 31
$ws.Range("K31").Value = 19
$ws.Range("L31").Value = 22.25

# Row 32
$ws.Range("K32").Value = 22
$ws.Range("L32").Value = 57.42

# Row 35 (totals)
$ws.Range("K35").Value = 356
$ws.Range("L35").Value = 449.71
